# The deck ships two themes:
#   ppt/theme/theme1.xml  -> "Luxe"    (the live theme, used by the slide master)
#   ppt/theme/theme2.xml  -> "Default" (an unused theme, only referenced by the notes master)
#
# The authoritative edit swaps the two themes' contents: theme1.xml becomes the
# "Default" colour scheme and theme2.xml becomes the "Luxe" colour scheme (font
# scheme and format scheme are identical between the two themes already, so the
# only real difference is the 12 colour-scheme slots).
#
# PowerPoint's automation surface only exposes a single, editable theme object
# (reached via the slide master / notes master / handout master / any slide) -
# it is backed by ppt/theme/theme1.xml, the deck's live/rendered theme. So we
# repoint that live theme at the "Default" palette, which is the half of the
# swap that actually affects what the presentation looks like.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Helper: RGB() builds a PowerPoint BGR-packed long from R,G,B bytes, matching
# what ThemeColorScheme.Colors(i).RGB expects/returns.
function Hex-RGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Default" theme colour scheme (was previously theme2.xml's palette),
# in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$defaultColors = @(
    "000000",
    "FFFFFF",
    "158158",
    "F3F3F3",
    "058DC7",
    "50B432",
    "ED561B",
    "EDEF00",
    "24CBE5",
    "64E572",
    "2200CC",
    "551A8B"
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = Hex-RGB $defaultColors[$i - 1]
}
